# "add tset case 6" - Guru99 Live Project Demo TestCase workbook
# Adds two new test cases (rows 7 & 8) below the existing "share wishlist" case (row 6),
# tweaks row 6 so its Steps column reflects logging in (instead of re-registering),
# and flags row 6 as done via a new column F checkmark.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New Case 6 (row 7): "Verify user is able to purchase ... " ---
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "Verify user is able to purchase product using registed email id(Use chrome browser)"

# --- Row 6 (Case 5): Steps column now assumes an existing login instead of fresh registration ---
$ws.Cells.Item(6, 3).Value = "1. Goto http://live.demoguru99.com`n2. Click on my account link`n4. login`n5. goto Tv menu`n6. Add product in your wish list`n7. Click share wishlist`n8. In next page enter Email and a message and click share wishlist`n9. Check wishlist is shared"
$ws.Cells.Item(6, 6).Value = 1
$ws.Rows.Item(6).RowHeight = 160

# --- finish New Case 6 (row 7) ---
$ws.Cells.Item(7, 3).Value = "1. Goto http://live.demoguru99.com`n2. Click on my account link`n3. Loing in application`n4. Click on my wishlist link`n5. In next oage, click add to cart link`n6. Click proceed to checkout `n7.enter shipping information`n8. click estimate`n9. verify shipping cost generated`n10. select shipping cost.update total`n11. verify shipping cost is add t total`n12. click 'Proceed to checkout`"`n13. Enter billing information`n14. In shipping method, click continue`n15. In payment informtaion select 'Check/Money order' radio button. Click continue`n16. Click 'Place order' button`n17. Verify order is generated. Note the order number"
$ws.Cells.Item(7, 4).Value = "1. Shipping information country = united states `nstate =  new york`nzip = 542896`naddress = ABC `ncity = new york`nTelephone = 12345678"
$ws.Cells.Item(7, 5).Value = "1. flat rate shipping of `$5 is generated`n2.Shipping cost is added to total product cost`n3.Order is placed. Order number is generated"
$ws.Rows.Item(7).RowHeight = 351

# --- New Case 7 (row 8): "Verify that you will be able ot save previously placed order as a pdf file" ---
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "Verify that you will be able ot save previously placed order as a pdf file "
$ws.Cells.Item(8, 3).Value = "1. Go to http://live. Demoguru99.com`n2. Click on my account link `n3. Login in application`n4. Click on 'My Orders'`n5. Click on 'View order'`n6. Verify the previously created order is displayed in 'recent orders' table and status is pending`n7. Click on 'print order' link`n8. Verify order is saved as PDF"
$ws.Cells.Item(8, 4).Value = "Use FireFox"

# --- Row 9: E9 carries a quote-prefixed (empty) text entry ---
$ws.Cells.Item(9, 5).Value = "'"
$ws.Rows.Item(9).RowHeight = 20

# --- finish New Case 7 (row 8) Expect column ---
$ws.Cells.Item(8, 5).Value = "1. Previously created order is displayed in 'Recent orders' table and status is pending `n2. Order is saved as PDF"
$ws.Rows.Item(8).RowHeight = 180

# --- View state: move viewport/selection to the newly added rows ---
$ws.Range("E9").Select()

